$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 13.4877
$ws.Range("E6").Value = 12.37330000000001
$ws.Range("E7").Value = 11.917
$ws.Range("E16").Value = 11.96240000000001
$ws.Range("E20").Value = 13.33999999999999
